$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns.
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (e.g. "1.000", "0.4500") keep their exact
# textual representation instead of being coerced to numbers, and
# Style is reset to "Normal" afterwards so no stray formatting is left
# behind on cells whose value does not otherwise change.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '30.298.81'
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -3.49%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.934.89'
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -3.42%  '
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.31%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '248.05'
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -2.95%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.7215'
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -11.39%  '
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3294'
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -8.12%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '26.71'
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +3.28%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.06844'
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -2.45%  '
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.8099'
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -4.50%  '
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07977'
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.932.95'
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.456'
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -1.49%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '95.04'
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -6.49%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '14.62'
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +4.14%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '263.79'
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '30.290.08'
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -3.47%  '
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.847'
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '2.187.03'
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -2.77%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.9991'
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.9993'
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +0.33%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.940'
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -1.20%  '
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.765'
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -1.02%  '
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '160.39'
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -2.26%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.360'
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +4.61%  '
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.1348'
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -11.84%  '
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.11'
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -5.33%  '
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +0.78%  '
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.557'
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -1.21%  '
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.429'
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -4.05%  '
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.234'
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -3.00%  '
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.05105'
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -2.10%  '
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.213'
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.77%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7488'
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -2.03%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.729'
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -0.98%  '
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01951'
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -3.22%  '
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.824'
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -2.99%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '81.24'
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +3.08%  '
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.592'
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -0.98%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.4500'
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -5.87%  '
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.025'
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -5.28%  '
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.8392'
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -2.75%  '
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '102.62'
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -1.98%  '
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.763'
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -2.17%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.371'
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -2.45%  '
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '36.42'
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -1.53%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.504'
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +2.92%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.4140'
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -6.21%  '
$c.Style = "Normal"

